$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap columns A ("FICHA") and B ("NOMBRE") -- including the narrow
# custom column width that used to live on column A -- by inserting a
# blank column in front of A (which shifts the old A, width and all,
# into B), copying the old B value into the new A, and then collapsing
# the duplicate (old B, now shifted to C) back out so every column from
# C onward lines back up with the original layout.
$oldB = $ws.Range("B1").Value()

$ws.Columns.Item(1).Insert()
$ws.Range("A1").Value = $oldB
$ws.Columns.Item(3).Delete()

# Restore the selection/active cell as recorded in the sheet view
$ws.Range("C4").Select()
